# Model_Data_Base_methane.xlsx fix-up
#
# "Error in methane/jet_fuel input + temporal resolution update"
#
# The "Connections" sheet's Table13 had an accidental extra column
# ("Column1") that was holding the "initial_connections_invested_available"
# flag values (1), while the real "initial_connections_invested_available"
# column was instead holding the "connection_investment_tech_lifetime"
# values ("40Y"/"50Y"). This script moves the data back into the correct
# columns, fills in the one row that was missing its
# "initial_connections_invested_available" value, and removes the stray
# "Column1" column from the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Connections")

# Column Y = connection_investment_tech_lifetime
# Column Z = initial_connections_invested_available
# Column AA = stray "Column1" (to be removed)

# Rows 2, 3 and 5 had their "connection_investment_tech_lifetime" value
# ("40Y"/"50Y") mistakenly entered one column to the right, in
# "initial_connections_invested_available" (Z). Move it back to Y.
$ws.Range("Y2").Value2 = $ws.Range("Z2").Value2
$ws.Range("Y3").Value2 = $ws.Range("Z3").Value2
$ws.Range("Y5").Value2 = $ws.Range("Z5").Value2

# The real "initial_connections_invested_available" values (1) had ended up
# in the stray "Column1" (AA). Move them back to Z.
$ws.Range("Z2").Value2 = $ws.Range("AA2").Value2
$ws.Range("Z3").Value2 = $ws.Range("AA3").Value2
$ws.Range("Z5").Value2 = $ws.Range("AA5").Value2

# Row 4 (pl_ch4_st) was missing its "initial_connections_invested_available"
# value entirely - fill it in.
$ws.Range("Z4").Value2 = 1

# Remove the stray extra "Column1" table column now that its data has been
# relocated.
$lo = $ws.ListObjects.Item(1)
$col = $lo.ListColumns.Item("Column1")
$col.Delete()

# Reflect that "Connections" is now the sheet the user was last working on.
$ws1 = $wb.Worksheets.Item("Units")
$ws1.Activate() | Out-Null
$ws1.Range("G11").Select() | Out-Null

$ws.Activate() | Out-Null
$ws.Range("Z5").Select() | Out-Null
